$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Day 1 Assignments")
$ws.Activate()

$ws.Range("D29").Value = "1. NameNotFound"
$ws.Range("D30").Value = "2. LOcationNotFoundException"
$ws.Range("D31").Value = "3. INvalidAmountException"

$ws.Range("D33").Value = "CRUD :"
$ws.Range("D34").Value = "create"
$ws.Range("D35").Value = "read"
$ws.Range("D36").Value = "update"
$ws.Range("D37").Value = "delete"

$ws.Range("D38").Value = "CRUD : Subscription : store the data in AL"
$ws.Range("D39").Value = "CRUD : Movies : store the data in HS"
$ws.Range("D40").Value = "CRUD : Series : Store the data in TS"

$ws.Range("D41").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
